$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 298.33334
$ws.Range("I53").Value = 288.75
$ws.Range("K53").Value = 288.75
$ws.Range("M53").Value = 348.25
# Row 74
$ws.Range("H74").Value = 13160
$ws.Range("I74").Value = 13160
$ws.Range("K74").Value = 13160
$ws.Range("M74").Value = -12224
# Row 76
$ws.Range("H76").Value = 4003
$ws.Range("I76").Value = 4003
$ws.Range("K76").Value = 4003
$ws.Range("M76").Value = -3688
# Row 77
$ws.Range("H77").Value = 13160
$ws.Range("I77").Value = 13160
$ws.Range("K77").Value = 65800
$ws.Range("M77").Value = -61120
# Row 79
$ws.Range("H79").Value = 4003
$ws.Range("I79").Value = 4003
$ws.Range("K79").Value = 4003
$ws.Range("M79").Value = -2911
# Row 80
$ws.Range("H80").Value = 4799.2
$ws.Range("I80").Value = 2999
$ws.Range("J80").Value = 5249.25
$ws.Range("K80").Value = 8997
$ws.Range("L80").Value = 15747.75
$ws.Range("M80").Value = -7999
$ws.Range("N80").Value = -17743.75
# Row 83
$ws.Range("H83").Value = 4799.2
$ws.Range("I83").Value = 2999
$ws.Range("J83").Value = 5249.25
$ws.Range("K83").Value = 26991
$ws.Range("L83").Value = 47243.25
$ws.Range("M83").Value = -21999
$ws.Range("N83").Value = -57227.25
# Row 98
$ws.Range("H98").Value = 2487.8333
$ws.Range("I98").Value = 856.6
$ws.Range("K98").Value = 856.6
$ws.Range("M98").Value = 641.4
# Row 107
$ws.Range("H107").Value = 1537.0555
$ws.Range("I107").Value = 1509.8235
$ws.Range("K107").Value = 1509.8235
$ws.Range("M107").Value = 410.1765
# Row 113
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746
# Row 122
$ws.Range("H122").Value = 2487.8333
$ws.Range("I122").Value = 856.6
$ws.Range("K122").Value = 2569.8
$ws.Range("M122").Value = -119.8000000000002
# Row 137
$ws.Range("H137").Value = 1767.3529
$ws.Range("I137").Value = 1553.2858
$ws.Range("J137").Value = 2766.3333
$ws.Range("K137").Value = 4659.857400000001
$ws.Range("L137").Value = 8298.999899999999
$ws.Range("M137").Value = -2109.857400000001
$ws.Range("N137").Value = -13398.9999
# Row 138
$ws.Range("H138").Value = 7098.8887
$ws.Range("J138").Value = 7402.9165
$ws.Range("L138").Value = 22208.7495
$ws.Range("N138").Value = -32488.7495

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7921.45
$ws.Range("I32").Value = 7921.45
$ws.Range("K32").Value = 7921.45
$ws.Range("M32").Value = -7634.45
# Row 46
$ws.Range("H46").Value = 13855.6
$ws.Range("I46").Value = 10069
$ws.Range("J46").Value = 16380
$ws.Range("K46").Value = 10069
$ws.Range("L46").Value = 16380
$ws.Range("M46").Value = -9750
$ws.Range("N46").Value = -17018
# Row 63
$ws.Range("H63").Value = 7298.8887
$ws.Range("I63").Value = 947.3333
$ws.Range("J63").Value = 20002
$ws.Range("K63").Value = 947.3333
$ws.Range("L63").Value = 20002
$ws.Range("M63").Value = -261.3333
$ws.Range("N63").Value = -21374
# Row 66
$ws.Range("H66").Value = 7298.8887
$ws.Range("I66").Value = 947.3333
$ws.Range("J66").Value = 20002
$ws.Range("K66").Value = 4736.6665
$ws.Range("L66").Value = 100010
$ws.Range("M66").Value = -1304.6665
$ws.Range("N66").Value = -106874
# Row 97
$ws.Range("H97").Value = 4145.5713
$ws.Range("J97").Value = 13750
$ws.Range("L97").Value = 13750
$ws.Range("N97").Value = -14742

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 60074
$ws.Range("J35").Value = 60074
$ws.Range("L35").Value = 60074
$ws.Range("N35").Value = -60694
# Row 64
$ws.Range("H64").Value = 3722.6667
$ws.Range("J64").Value = 5482.6665
$ws.Range("L64").Value = 5482.6665
$ws.Range("N64").Value = -5932.6665
# Row 67
$ws.Range("H67").Value = 3722.6667
$ws.Range("J67").Value = 5482.6665
$ws.Range("L67").Value = 5482.6665
$ws.Range("N67").Value = -7042.6665
# Row 94
$ws.Range("H94").Value = 4078.5386
$ws.Range("I94").Value = 4184.636
$ws.Range("K94").Value = 4184.636
$ws.Range("M94").Value = -3733.636
# Row 107
$ws.Range("H107").Value = 1996.6666
$ws.Range("I107").Value = 1947.5
$ws.Range("K107").Value = 1947.5
$ws.Range("M107").Value = -27.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1351.5454
$ws.Range("I107").Value = 886.8
$ws.Range("K107").Value = 886.8
$ws.Range("M107").Value = 1033.2
# Row 141
$ws.Range("H141").Value = 598998.5
$ws.Range("I141").Value = 300000
$ws.Range("J141").Value = 698664.7
$ws.Range("K141").Value = 300000
$ws.Range("L141").Value = 698664.7
$ws.Range("M141").Value = -294820
$ws.Range("N141").Value = -709024.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 97.083336
$ws.Range("I2").Value = 86
$ws.Range("J2").Value = 102.625
$ws.Range("K2").Value = 516
$ws.Range("L2").Value = 615.75
$ws.Range("M2").Value = -403
$ws.Range("N2").Value = -841.75
# Row 34
$ws.Range("H34").Value = 7499.8335
$ws.Range("J34").Value = 11000
$ws.Range("L34").Value = 33000
$ws.Range("N34").Value = -33168
# Row 97
$ws.Range("H97").Value = 1232.7142
$ws.Range("J97").Value = 840.75
$ws.Range("L97").Value = 2522.25
$ws.Range("N97").Value = -3514.25
# Row 98
$ws.Range("H98").Value = 3292.3333
$ws.Range("I98").Value = 3444.75
$ws.Range("J98").Value = 2987.5
$ws.Range("K98").Value = 10334.25
$ws.Range("L98").Value = 8962.5
$ws.Range("M98").Value = -8836.25
$ws.Range("N98").Value = -11958.5
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
# Row 139
$ws.Range("H139").Value = 1940
$ws.Range("I139").Value = 1940
$ws.Range("K139").Value = 5820
$ws.Range("M139").Value = -680

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 6141.25
$ws.Range("I102").Value = 6141.25
$ws.Range("K102").Value = 6141.25
$ws.Range("M102").Value = -4519.25
# Row 107
$ws.Range("H107").Value = 5331.6665
$ws.Range("I107").Value = 7500
$ws.Range("J107").Value = 995
$ws.Range("K107").Value = 7500
$ws.Range("L107").Value = 995
$ws.Range("M107").Value = -5580
$ws.Range("N107").Value = -4835

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2187.5
$ws.Range("I16").Value = 1653.1428
$ws.Range("J16").Value = 3434.3333
$ws.Range("K16").Value = 1653.1428
$ws.Range("L16").Value = 3434.3333
$ws.Range("M16").Value = -1483.1428
$ws.Range("N16").Value = -3774.3333
# Row 22
$ws.Range("H22").Value = 4771.722
$ws.Range("I22").Value = 2446.818
$ws.Range("J22").Value = 8425.143
$ws.Range("K22").Value = 2446.818
$ws.Range("L22").Value = 8425.143
$ws.Range("M22").Value = -2151.818
$ws.Range("N22").Value = -9015.143
# Row 27
$ws.Range("H27").Value = 4771.722
$ws.Range("I27").Value = 2446.818
$ws.Range("J27").Value = 8425.143
$ws.Range("K27").Value = 2446.818
$ws.Range("L27").Value = 8425.143
$ws.Range("M27").Value = -2339.818
$ws.Range("N27").Value = -8639.143
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 132
$ws.Range("H132").Value = 5999
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# Row 136
$ws.Range("H136").Value = 4094.3
$ws.Range("I136").Value = 1811.75
$ws.Range("K136").Value = 5435.25
$ws.Range("M136").Value = -2885.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1171.8
$ws.Range("I122").Value = 1171.8
$ws.Range("K122").Value = 3515.4
$ws.Range("M122").Value = -1065.4
